$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z7000").Value = "0.51279"
$ws.Range("Z7000").NumberFormat = "General"
